# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'34.776.67"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.853.12"
$ws.Range("E3").Value = "  +2.29%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +0.79%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.74%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.26%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'41.92"
$ws.Range("E8").Value = "  +14.62%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +3.82%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.08%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +3.59%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "'2.120.06"
$ws.Range("E12").Value = "  +2.27%  "

# Row 13 - was WrappedEther, now Chainlink (rows 13/14 swapped content)
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.39"
$ws.Range("E13").Value = "  +0.57%  "

# Row 14 - was Chainlink, now WrappedEther
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.849.13"
$ws.Range("E14").Value = "  +1.71%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "'4.74"
$ws.Range("E15").Value = "  +6.64%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "'0.660"
$ws.Range("E16").Value = "  +4.45%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'34.723.71"
$ws.Range("E17").Value = "  +0.55%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'69.03"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'244.93"
$ws.Range("E19").Value = "  +0.66%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.0₃0789"
$ws.Range("E20").Value = "  +1.42%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'12.16"
$ws.Range("E21").Value = "  +8.06%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.77"
$ws.Range("E22").Value = "  +16.02%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.32%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  -1.47%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'172.56"
$ws.Range("E25").Value = "  +0.54%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'7.93"
$ws.Range("E26").Value = "  +0.86%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'17.85"
$ws.Range("E27").Value = "  +3.11%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.06%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.34%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'1.34"
$ws.Range("E30").Value = "  +8.31%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  +2.86%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'4.01"
$ws.Range("E32").Value = "  +2.28%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0533"
$ws.Range("E33").Value = "  +3.26%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "  +5.43%  "

# Row 35 - Aave
$ws.Range("D35").Value = "'91.47"
$ws.Range("E35").Value = "  +12.44%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.670"
$ws.Range("E36").Value = "  +2.51%  "

# Row 37 - was Maker, now ARBITRUM (rows 37/38 swapped content)
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = "  +10.28%  "

# Row 38 - was ARBITRUM, now Maker
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'1.345.33"
$ws.Range("E38").Value = "  -1.37%  "

# Row 39 - was RenderToken, now TrustWalletToken (rows 39/40 swapped content)
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.08"
$ws.Range("E39").Value = "  +1.67%  "

# Row 40 - was TrustWalletToken, now RenderToken
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.44"
$ws.Range("E40").Value = "  +3.02%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "'0.0194"
$ws.Range("E41").Value = "  +3.78%  "

# Row 42 - InjectiveProtocol
$ws.Range("D42").Value = "'14.82"
$ws.Range("E42").Value = "  +8.91%  "

# Row 43 - WEMIXToken
$ws.Range("E43").Value = "  +7.47%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  +2.29%  "

# Row 45 - HuobiToken
$ws.Range("E45").Value = "  +0.85%  "

# Row 46 - Kaspa
$ws.Range("D46").Value = "'0.0518"
$ws.Range("E46").Value = "  +3.49%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "'2.017.71"
$ws.Range("E47").Value = "  +2.27%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "'6.04"
$ws.Range("E48").Value = "  +3.72%  "

# Row 49 - PaxDollar
$ws.Range("E49").Value = "  +0.28%  "

# Row 50 - Quant
$ws.Range("D50").Value = "'102.66"
$ws.Range("E50").Value = "  -0.11%  "

# Row 51 - Aptos
$ws.Range("D51").Value = "'7.30"
$ws.Range("E51").Value = "  +5.56%  "
